$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.666.90"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "2.646.36"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.48"
$ws.Range("E5").Value = "  +1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.86"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  -1.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.48"
$ws.Range("E9").Value = "  -3.13%  "

$ws.Range("E10").Value = "  +3.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.349"
$ws.Range("E11").Value = "  +0.54%  "

$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").Value = "3.108.65"
$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("D14").Value = "60.683.59"
$ws.Range("E14").Value = "  -0.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.85"
$ws.Range("E15").Value = "  +0.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000142"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").Value = "2.638.20"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.74"
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "351.78"
$ws.Range("E19").Value = "  -1.67%  "

$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.26"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").Value = "0.0₃0848"
$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.24"
$ws.Range("E28").Value = "  -2.11%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.15"
$ws.Range("E30").Value = "  +3.78%  "

$ws.Range("E31").Value = "  +1.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.34"
$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.79"
$ws.Range("E33").Value = "  -1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.09"
$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("E35").Value = "  -1.13%  "

$ws.Range("E36").Value = "  +5.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.887"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.84"
$ws.Range("E38").Value = "  +1.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "305.67"
$ws.Range("E39").Value = "  +4.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.75"
$ws.Range("E40").Value = "  -0.23%  "

$ws.Range("E41").Value = "  -2.06%  "

$ws.Range("E42").Value = "  +1.96%  "

$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.07"
$ws.Range("E44").Value = "  +1.32%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0558"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("E47").Value = "  +1.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.88"
$ws.Range("E48").Value = "  -2.34%  "

$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.06"
$ws.Range("E50").Value = "  +1.95%  "

$ws.Range("D51").Value = "1.983.41"
$ws.Range("E51").Value = "  -0.71%  "
